# Update date and division problems/answers to the new set.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

Replace-Text "2024-11-12 Tuesday" "2024-11-13 Wednesday"
Replace-Text "438÷5=87, 3" "464÷2=232, 0"
Replace-Text "671÷8=83, 7" "363÷7=51, 6"
Replace-Text "996÷8=124, 4" "849÷3=283, 0"
Replace-Text "647÷3=215, 2" "506÷8=63, 2"
Replace-Text "319÷7=45, 4" "533÷8=66, 5"
Replace-Text "164÷5=32, 4" "203÷5=40, 3"
Replace-Text "301÷8=37, 5" "285÷7=40, 5"
Replace-Text "888÷8=111, 0" "146÷3=48, 2"
Replace-Text "835÷7=119, 2" "296÷8=37, 0"
Replace-Text "158÷9=17, 5" "116÷3=38, 2"
Replace-Text "680÷3=226, 2" "644÷6=107, 2"
Replace-Text "923÷8=115, 3" "810÷3=270, 0"
Replace-Text "449÷5=89, 4" "111÷5=22, 1"
Replace-Text "771÷8=96, 3" "767÷3=255, 2"
Replace-Text "336÷4=84, 0" "393÷2=196, 1"
Replace-Text "850÷7=121, 3" "381÷9=42, 3"
Replace-Text "825÷6=137, 3" "495÷9=55, 0"
Replace-Text "558÷2=279, 0" "819÷2=409, 1"
Replace-Text "512÷6=85, 2" "634÷4=158, 2"
Replace-Text "206÷5=41, 1" "912÷4=228, 0"
Replace-Text "313÷3=104, 1" "258÷6=43, 0"
Replace-Text "405÷6=67, 3" "653÷2=326, 1"
Replace-Text "677÷9=75, 2" "659÷4=164, 3"
Replace-Text "204÷9=22, 6" "929÷2=464, 1"
Replace-Text "943÷2=471, 1" "115÷7=16, 3"
